# Update column C values (Name of Algo data) for specific rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8   = -12.741
    10  = -13.294
    12  = -10.616
    18  = -11.865
    25  = -11.96
    37  = -12.863
    55  = -13.837
    68  = -11.134
    77  = -13.145
    78  = -12.892
    79  = -12.34
    80  = -12.858
    81  = -13.404
    82  = -11.895
    84  = -13.636
    101 = -12.37
    102 = -12.786
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
